# Insert a new weekly price record as row 154 in the "Acelga" dataset,
# pushing the existing rows 154:224 down to 155:225.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("154:154").Insert()

$ws.Range("A154").Value = 10
$ws.Range("B154").Value = "Vega Modelo de Temuco"
$ws.Range("C154").Value = "La Araucanía"
$ws.Range("D154").Value = 44523
$ws.Range("E154").Value = 9
$ws.Range("F154").Value = 100112009
$ws.Range("G154").Value = "Acelga"
$ws.Range("H154").Value = "Sin especificar"
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 50
$ws.Range("K154").Value = 9000
$ws.Range("L154").Value = 9000
$ws.Range("M154").Value = 9000
$ws.Range("N154").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O154").Value = "Provincia de Cautín"
$ws.Range("P154").Value = 750
$ws.Range("Q154").Value = 12
$ws.Range("R154").Value = "Hortaliza"
